$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 previously held the text "R40"; it is now replaced with the text "1".
# Prefixing with an apostrophe forces Excel to store it as text (matching the
# source shared-string entry) instead of silently coercing it to the number 1.
$ws.Range("B11").Value = "'1"
